# Natmi following Dr Hou advice:
# The NATMI LR-pair output for Il4-Il4ra was recomputed with an extra
# sending/target cluster ("M2") added to the analysis. This rewrites every
# data row (2-13) with the recalculated statistics and appends the new
# rows (14-17) needed for the "M2"-as-sender combinations against every
# target cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Il4", "Il4ra", "ECs", 3, 1, 0.9334919999999999, 2.800476, 0.209095668602663, 0.209095668602663, 3, 1, 43.23621900000001, 129.708657, 0.4133516515244107, 0.4133516515244106, 40.360664546748, 363.245980920732, 0.08643003994351162, 0.08643003994351159),
  @("ECs", "Il4", "Il4ra", "FAPs", 3, 1, 0.9334919999999999, 2.800476, 0.209095668602663, 0.209095668602663, 3, 1, 21.35406233333333, 64.062187, 0.2041514530268833, 0.2041514530268833, 19.933846355668, 179.404617201012, 0.04268718456686132, 0.04268718456686131),
  @("ECs", "Il4", "Il4ra", "M2", 3, 1, 0.9334919999999999, 2.800476, 0.209095668602663, 0.209095668602663, 3, 1, 32.81697399999999, 98.45092199999999, 0.3137404406461546, 0.3137404406461545, 30.63438269320799, 275.7094442388719, 0.0656017672046018, 0.06560176720460179),
  @("ECs", "Il4", "Il4ra", "sCs", 3, 1, 0.9334919999999999, 2.800476, 0.209095668602663, 0.209095668602663, 3, 1, 7.191864666666667, 21.575594, 0.06875645480255156, 0.06875645480255155, 6.713548131415999, 60.42193318274399, 0.0143766768876883, 0.01437667688768829),
  @("FAPs", "Il4", "Il4ra", "ECs", 3, 1, 1.964334, 5.893002, 0.4399970552387631, 0.4399970552387631, 3, 1, 43.23621900000001, 129.708657, 0.4133516515244107, 0.4133516515244106, 84.930375013146, 764.3733751183141, 0.1818735094488201, 0.18187350944882),
  @("FAPs", "Il4", "Il4ra", "FAPs", 3, 1, 1.964334, 5.893002, 0.4399970552387631, 0.4399970552387631, 3, 1, 21.35406233333333, 64.062187, 0.2041514530268833, 0.2041514530268833, 41.946510679486, 377.518596115374, 0.08982603815454335, 0.08982603815454332),
  @("FAPs", "Il4", "Il4ra", "M2", 3, 1, 1.964334, 5.893002, 0.4399970552387631, 0.4399970552387631, 3, 1, 32.81697399999999, 98.45092199999999, 0.3137404406461546, 0.3137404406461545, 64.463497805316, 580.171480247844, 0.13804486999362, 0.1380448699936199),
  @("FAPs", "Il4", "Il4ra", "sCs", 3, 1, 1.964334, 5.893002, 0.4399970552387631, 0.4399970552387631, 3, 1, 7.191864666666667, 21.575594, 0.06875645480255156, 0.06875645480255155, 14.127224288132, 127.145018593188, 0.0302526376417798, 0.03025263764177979),
  @("M2", "Il4", "Il4ra", "ECs", 2, 0.6666666666666666, 0.1531806666666667, 0.459542, 0.03431139625585257, 0.03431139625585256, 3, 1, 43.23621900000001, 129.708657, 0.4133516515244107, 0.4133516515244106, 6.622952850566, 59.60657565509401, 0.01418267230846514, 0.01418267230846513),
  @("M2", "Il4", "Il4ra", "FAPs", 2, 0.6666666666666666, 0.1531806666666667, 0.459542, 0.03431139625585257, 0.03431139625585256, 3, 1, 21.35406233333333, 64.062187, 0.2041514530268833, 0.2041514530268833, 3.271029504261555, 29.439265538354, 0.007004721401013467, 0.007004721401013464),
  @("M2", "Il4", "Il4ra", "M2", 2, 0.6666666666666666, 0.1531806666666667, 0.459542, 0.03431139625585257, 0.03431139625585256, 3, 1, 32.81697399999999, 98.45092199999999, 0.3137404406461546, 0.3137404406461545, 5.026925955302666, 45.242333597724, 0.010764872580496, 0.010764872580496),
  @("M2", "Il4", "Il4ra", "sCs", 2, 0.6666666666666666, 0.1531806666666667, 0.459542, 0.03431139625585257, 0.03431139625585256, 3, 1, 7.191864666666667, 21.575594, 0.06875645480255156, 0.06875645480255155, 1.101654624216444, 9.914891617948, 0.002359129965877964, 0.002359129965877963),
  @("sCs", "Il4", "Il4ra", "ECs", 3, 1, 1.413418666666667, 4.240256, 0.3165958799027214, 0.3165958799027214, 3, 1, 43.23621900000001, 129.708657, 0.4133516515244107, 0.4133516515244106, 61.11087901068802, 549.9979110961922, 0.1308654298236139, 0.1308654298236138),
  @("sCs", "Il4", "Il4ra", "FAPs", 3, 1, 1.413418666666667, 4.240256, 0.3165958799027214, 0.3165958799027214, 3, 1, 21.35406233333333, 64.062187, 0.2041514530268833, 0.2041514530268833, 30.18223031109689, 271.640072799872, 0.06463350890446523, 0.06463350890446522),
  @("sCs", "Il4", "Il4ra", "M2", 3, 1, 1.413418666666667, 4.240256, 0.3165958799027214, 0.3165958799027214, 3, 1, 32.81697399999999, 98.45092199999999, 0.3137404406461546, 0.3137404406461545, 46.38412363511466, 417.457112716032, 0.09932893086743685, 0.09932893086743683),
  @("sCs", "Il4", "Il4ra", "sCs", 3, 1, 1.413418666666667, 4.240256, 0.3165958799027214, 0.3165958799027214, 3, 1, 7.191864666666667, 21.575594, 0.06875645480255156, 0.06875645480255155, 10.16511576800711, 91.48604191206401, 0.02176801030720551, 0.0217680103072055)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
